# Update instructions and hardware.
$wb = $excel.ActiveWorkbook

# Fix typo in the "Instructions FR" sheet: remove duplicated word "boitier"
# in the first step's label ("Impression 3D du boitier boitier" ->
# "Impression 3D du boitier").
$wsFR = $wb.Worksheets.Item("Instructions FR")
$wsFR.Range("A2").Value = "Impression 3D du boitier"

# Restore selection / active-sheet state to match the saved workbook: the
# French sheet becomes the active tab with B17 selected, the English sheet
# is left with B10 selected and is no longer the active tab.
$wsEN = $wb.Worksheets.Item("Instructions EN")
$wsEN.Range("B10").Select()

$wsFR.Activate()
$wsFR.Range("B17").Select()
